$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 199, shifting rows 199:324 down to 200:325.
$ws.Rows(199).Insert()

# Populate the newly inserted row 199 with the new weekly record.
$ws.Range("A199").Value = 3
$ws.Range("B199").Value = "Femacal de La Calera"
$ws.Range("C199").Value = "Coquimbo"
$ws.Range("D199").Value = 44582
$ws.Range("E199").Value = 5
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100108
$ws.Range("H199").Value = "Tropicales y subtropicales"
$ws.Range("I199").Value = 100108002
$ws.Range("J199").Value = "Mango"
$ws.Range("K199").Value = "Sin especificar"
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 160
$ws.Range("N199").Value = 7000
$ws.Range("O199").Value = 7000
$ws.Range("P199").Value = 7000
$ws.Range("Q199").Value = "`$/bandeja 4 kilos"
$ws.Range("R199").Value = "Perú"
$ws.Range("S199").Value = 1750
$ws.Range("T199").Value = 4
